# Refresh the cryptos table (prices in column D, 1h-volume deltas in column E)
# to match the latest scrape. Row 33/34 also swap content (Filecoin <-> WEMIXTOKEN)
# because the source ranking changed order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Note: some "Price" values are plain digit-and-dot strings (e.g. "1.001") that
# Excel would otherwise auto-convert to a Double when assigned via .Value. Prefixing
# with an apostrophe (as if the user typed '1.001 into the cell) keeps them as text,
# matching the original inline-string cell content.

$ws.Range('D2').Value = '22.423.05'

$ws.Range('D3').Value = '1.564.19'
$ws.Range('E3').Value = '  -0.50%  '

$ws.Range('D4').Value = '''1.001'
$ws.Range('E4').Value = '  -0.19%  '

$ws.Range('E5').Value = '  -0.10%  '

$ws.Range('D6').Value = '''285.75'
$ws.Range('E6').Value = '  -1.94%  '

$ws.Range('D7').Value = '''0.3623'
$ws.Range('E7').Value = '  -2.84%  '

$ws.Range('E8').Value = '  -2.67%  '

$ws.Range('D9').Value = '''0.3340'
$ws.Range('E9').Value = '  -1.49%  '

$ws.Range('D10').Value = '''1.124'
$ws.Range('E10').Value = '  -1.31%  '

$ws.Range('D11').Value = '''0.07384'
$ws.Range('E11').Value = '  -2.48%  '

$ws.Range('E12').Value = '  -0.20%  '

$ws.Range('D13').Value = '''20.72'
$ws.Range('E13').Value = '  -2.85%  '

$ws.Range('D14').Value = '''5.936'
$ws.Range('E14').Value = '  -0.88%  '

$ws.Range('D15').Value = '''6.888'
$ws.Range('E15').Value = '  -1.09%  '

$ws.Range('D16').Value = '1.562.85'
$ws.Range('E16').Value = '  -0.86%  '

$ws.Range('D17').Value = '''0.00001102'
$ws.Range('E17').Value = '  -1.53%  '

$ws.Range('D18').Value = '''88.24'
$ws.Range('E18').Value = '  -2.91%  '

$ws.Range('D19').Value = '''0.06701'
$ws.Range('E19').Value = '  -0.55%  '

$ws.Range('E20').Value = '  -0.18%  '

$ws.Range('D21').Value = '''6.334'
$ws.Range('E21').Value = '  +0.76%  '

$ws.Range('D22').Value = '''16.14'
$ws.Range('E22').Value = '  -1.20%  '

$ws.Range('E23').Value = '  -0.72%  '

$ws.Range('D24').Value = '22.419.08'
$ws.Range('E24').Value = '  -0.03%  '

$ws.Range('D25').Value = '''2.386'
$ws.Range('E25').Value = '  +2.23%  '

$ws.Range('D26').Value = '''2.548'
$ws.Range('E26').Value = '  -3.88%  '

$ws.Range('D27').Value = '''150.31'
$ws.Range('E27').Value = '  +1.39%  '

$ws.Range('D28').Value = '''19.31'
$ws.Range('E28').Value = '  -3.90%  '

$ws.Range('D29').Value = '''5.000'
$ws.Range('E29').Value = '  -0.19%  '

$ws.Range('D30').Value = '''123.45'
$ws.Range('E30').Value = '  -1.53%  '

$ws.Range('D31').Value = '1.739.43'
$ws.Range('E31').Value = '  -0.88%  '

$ws.Range('D32').Value = '''1.058'
$ws.Range('E32').Value = '  +1.22%  '

$ws.Range('B33').Value = 'WEMIXTOKEN'
$ws.Range('C33').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D33').Value = '''1.997'
$ws.Range('E33').Value = '  +1.36%  '

$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').Value = '''6.100'
$ws.Range('E34').Value = '  -0.83%  '

$ws.Range('D35').Value = '''9.745'
$ws.Range('E35').Value = '  -0.51%  '

$ws.Range('D36').Value = '''0.08282'
$ws.Range('E36').Value = '  -1.31%  '

$ws.Range('D37').Value = '''0.02398'
$ws.Range('E37').Value = '  -3.17%  '

$ws.Range('E38').Value = '  -2.68%  '

$ws.Range('D39').Value = '''0.06376'
$ws.Range('E39').Value = '  -2.16%  '

$ws.Range('D40').Value = '''1.296'
$ws.Range('E40').Value = '  -5.98%  '

$ws.Range('D41').Value = '''5.313'
$ws.Range('E41').Value = '  -2.69%  '

$ws.Range('D42').Value = '''11.11'
$ws.Range('E42').Value = '  -1.21%  '

$ws.Range('D43').Value = '''0.6071'
$ws.Range('E43').Value = '  -2.41%  '

$ws.Range('D44').Value = '''1.001'
$ws.Range('E44').Value = '  -0.36%  '

$ws.Range('D45').Value = '''13.81'
$ws.Range('E45').Value = '  -0.75%  '

$ws.Range('D46').Value = '''3.759'
$ws.Range('E46').Value = '  -1.39%  '

$ws.Range('D47').Value = '''0.5786'
$ws.Range('E47').Value = '  -0.09%  '

$ws.Range('D48').Value = '''2.012'
$ws.Range('E48').Value = '  -2.95%  '

$ws.Range('D49').Value = '''124.14'
$ws.Range('E49').Value = '  -4.12%  '

$ws.Range('D50').Value = '''1.217'
$ws.Range('E50').Value = '  +0.13%  '

$ws.Range('D51').Value = '''0.07198'
$ws.Range('E51').Value = '  -1.69%  '
